$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents (values/formulas) of A4:M11 but keep formatting/style
$ws.Range("A4:M11").ClearContents()

# Move the active selection to C4, matching the saved view state
$ws.Range("C4").Select()
